$wb = $excel.ActiveWorkbook

# Sheet "展览" (1st sheet): update F2:F5
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 419
$ws1.Range("F3").Value = 5209
$ws1.Range("F4").Value = 50
$ws1.Range("F5").Value = 54

# Sheet "全部类型" (4th sheet): update F2, F3, F5, F6
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 419
$ws4.Range("F3").Value = 5209
$ws4.Range("F5").Value = 50
$ws4.Range("F6").Value = 54
